$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.810.96"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.583.27"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'582.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").Value = "'146.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").Value = "'0.108"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "'27.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "3.046.99"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").Value = "62.766.58"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "2.584.41"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").Value = "'341.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "'4.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'67.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("D24").Value = "2.706.87"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("D26").Value = "'1.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'7.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "'8.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("D32").Value = "'466.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.91%  "
$ws.Range("D33").Value = "0.0₃0817"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").Value = "'176.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").Value = "'18.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").Value = "'4.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.23%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -3.13%  "
$ws.Range("D42").Value = "'158.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.33%  "
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").Value = "'0.630"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.28%  "
$ws.Range("D45").Value = "'20.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").Value = "'0.0236"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").Value = "'18.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'1.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'11.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.94%  "
